$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused "is_activated" column (E) entirely, shrinking the
# used range from A1:E3 down to A1:D3.
$ws.Range("E1:E3").Delete()

# "Number type input in access code info": level for the second data row
# should be the numeric 0, not 1.
$ws.Range("B3").Value = 0

# The price for the second data row ("test2") is no longer populated.
$ws.Range("C3").ClearContents()

# Leave the selection on the last touched cell.
$ws.Range("B3").Select()
